$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-08-19"

# Update the August row label to reflect the new "through" date
$ws.Range("A9").Value = "August (through 08-19)"

# Update August row (row 9) values for columns C:I (years 2016-2022)
$ws.Range("C9").Value = 45
$ws.Range("D9").Value = 48
$ws.Range("E9").Value = 31
$ws.Range("F9").Value = 27
$ws.Range("G9").Value = 122
$ws.Range("H9").Value = 99
$ws.Range("I9").Value = 101

# Update Total row (row 10) values for columns C:I (years 2016-2022)
$ws.Range("C10").Value = 347
$ws.Range("D10").Value = 513
$ws.Range("E10").Value = 456
$ws.Range("F10").Value = 331
$ws.Range("G10").Value = 743
$ws.Range("H10").Value = 1009
$ws.Range("I10").Value = 1072
